$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    # Force literal text storage so numeric-looking strings (e.g. "213.23")
    # are not auto-converted to numbers by the input parser.
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '27.353.95'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '1.656.66'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  -0.18%  '
Set-TextValue 'D5' '213.23'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  -0.18%  '
Set-TextValue 'D8' '23.69'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  -1.11%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('D12').Value = '1.890.39'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '1.656.16'
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.573'
$ws.Range('E14').Value = '  +3.77%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '4.08'
$ws.Range('E15').Value = '  -1.67%  '
Set-TextValue 'D16' '65.84'
$ws.Range('D17').Value = '27.354.01'
$ws.Range('E17').Value = '  -1.79%  '
Set-TextValue 'D18' '232.00'
$ws.Range('E18').Value = '  -7.61%  '
$ws.Range('E19').Value = '  -0.65%  '
Set-TextValue 'D20' '7.49'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('E21').Value = '  -0.15%  '
Set-TextValue 'D22' '4.38'
$ws.Range('E22').Value = '  -2.07%  '
Set-TextValue 'D23' '9.39'
$ws.Range('E23').Value = '  +0.43%  '
$ws.Range('E24').Value = '  -1.83%  '
Set-TextValue 'D25' '147.10'
$ws.Range('E25').Value = '  +0.15%  '
Set-TextValue 'D26' '7.16'
$ws.Range('E26').Value = '  -0.91%  '
Set-TextValue 'D27' '15.91'
$ws.Range('E27').Value = '  -2.52%  '
Set-TextValue 'D28' '0.999'
$ws.Range('E28').Value = '  -0.19%  '
Set-TextValue 'D29' '0.112'
$ws.Range('E29').Value = '  -0.45%  '
Set-TextValue 'D30' '0.0497'
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('E31').Value = '  -3.69%  '
Set-TextValue 'D32' '3.31'
$ws.Range('D33').Value = '1.452.29'
$ws.Range('E33').Value = '  +1.78%  '
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('E35').Value = '  +0.35%  '
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('E38').Value = '  -1.90%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  -0.18%  '
Set-TextValue 'D43' '65.28'
$ws.Range('E43').Value = '  -6.41%  '
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').Value = '1.798.00'
$ws.Range('E45').Value = '  -0.72%  '
Set-TextValue 'D46' '0.787'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('E47').Value = '  -0.56%  '
Set-TextValue 'D48' '88.19'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('E50').Value = '  -0.36%  '
Set-TextValue 'D51' '7.76'
$ws.Range('E51').Value = '  -0.98%  '
